$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: URL EPOS (C2) value stays the same text, shared-string index
#     just shifts because unused strings were pruned elsewhere; no actual
#     change needed to the cell content itself.

# --- Row 8 (header row): add headers for the new columns E..I
$ws.Range("E8").Value = "msisdn"
$ws.Range("F8").Value = "msi"
$ws.Range("G8").Value = "nip"
$ws.Range("H8").Value = "client"
$ws.Range("I8").Value = "portId"

# --- Row 9: update column B and add new columns E..I
$ws.Range("B9").Value = "501506363"
$ws.Range("E9").Value = "3045981684"
$ws.Range("F9").Value = "732111193278813"
$ws.Range("G9").Value = "81684"
$ws.Range("H9").Value = "1061520830"
$ws.Range("I9").Value = "00002201108240181684"

# --- Row 10: update column B and add new columns E..I
$ws.Range("B10").Value = "933727137"
$ws.Range("E10").Value = "3045981670"
$ws.Range("F10").Value = "732111193278811"
$ws.Range("G10").Value = "81670"
$ws.Range("H10").Value = "111295346"
$ws.Range("I10").Value = "00002201108240181670"

# --- Row 11: update column B only
$ws.Range("B11").Value = "333489166"

# --- Row 12: update column B only
$ws.Range("B12").Value = "725586919"

# --- Row 13: update column B only
$ws.Range("B13").Value = "25620076"

# --- Column widths: column E narrower (no longer bestFit), column I gets
#     a dedicated width; columns F..H revert to the sheet default width.
$ws.Columns.Item(5).ColumnWidth = 21.3
$ws.Columns.Item(9).ColumnWidth = 21.5

# --- View state: move the active selection to H18 (matches the saved
#     view in the diff).
$ws.Range("H18").Select()
